# Merge co-cooking rows with the spaces:
# Delete the co_cooking_A, co_cooking_B, co_cooking_C rows (rows 3, 5, 7)
# so the remaining rows shift up, matching the "merged" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

$ws.Rows.Item(5).Select()
